# Updated cryptos list on Mon Feb 19 06:59:13 UTC 2024 with GitHub Actions
#
# Applies the per-row Price (D) / Volume(1h) (E) updates, plus the
# Uniswap <-> ImmutableX row swap (rows 20-21: Coin name, Link and
# Price/Volume all move together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a "Price" (D) cell while
# keeping it a plain text cell (no residual cell style / number format),
# matching the original inlineStr text cells exactly.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Row 2: Bitcoin ---
Set-TextValue $ws.Range("D2") "52.428.38"
$ws.Range("E2").Value = "  +1.83%  "

# --- Row 3: Ethereum ---
Set-TextValue $ws.Range("D3") "2.922.45"
$ws.Range("E3").Value = "  +4.89%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.11%  "

# --- Row 5: BNB ---
Set-TextValue $ws.Range("D5") "352.01"
$ws.Range("E5").Value = "  -0.36%  "

# --- Row 6: Solana ---
Set-TextValue $ws.Range("D6") "112.50"
$ws.Range("E6").Value = "  +3.72%  "

# --- Row 7: XRP ---
Set-TextValue $ws.Range("D7") "0.562"
$ws.Range("E7").Value = "  +1.92%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  +0.11%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  +0.83%  "

# --- Row 10: Avalanche ---
Set-TextValue $ws.Range("D10") "40.16"
$ws.Range("E10").Value = "  +0.98%  "

# --- Row 11: Dogecoin ---
$ws.Range("E11").Value = "  +3.65%  "

# --- Row 12: TRON ---
$ws.Range("E12").Value = "  +0.58%  "

# --- Row 13: Chainlink ---
Set-TextValue $ws.Range("D13") "20.18"
$ws.Range("E13").Value = "  +1.30%  "

# --- Row 14: Polkadot ---
Set-TextValue $ws.Range("D14") "7.83"
$ws.Range("E14").Value = "  +1.48%  "

# --- Row 15: Wrapped liquid staked Ether 2.0 ---
Set-TextValue $ws.Range("D15") "3.381.86"
$ws.Range("E15").Value = "  +5.03%  "

# --- Row 16: Polygon ---
Set-TextValue $ws.Range("D16") "0.994"
$ws.Range("E16").Value = "  +5.74%  "

# --- Row 17: Wrapped Ether ---
Set-TextValue $ws.Range("D17") "2.906.86"
$ws.Range("E17").Value = "  +3.88%  "

# --- Row 18: Wrapped BTC ---
Set-TextValue $ws.Range("D18") "52.458.36"
$ws.Range("E18").Value = "  +1.97%  "

# --- Row 19: Internet Computer (DFINITY) ---
Set-TextValue $ws.Range("D19") "14.73"
$ws.Range("E19").Value = "  +10.00%  "

# --- Row 20/21: Uniswap and ImmutableX swap places ---
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D20") "7.71"
$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D21") "3.34"
$ws.Range("E21").Value = "  +5.41%  "

# --- Row 22: Shiba Inu ---
Set-TextValue $ws.Range("D22") "0.0₃0983"
$ws.Range("E22").Value = "  +1.37%  "

# --- Row 23: Litecoin ---
$ws.Range("E23").Value = "  +1.32%  "

# --- Row 24: Bitcoin Cash ---
Set-TextValue $ws.Range("D24") "271.62"
$ws.Range("E24").Value = "  +1.75%  "

# --- Row 25: PancakeSwap ---
Set-TextValue $ws.Range("D25") "2.79"
$ws.Range("E25").Value = "  +1.76%  "

# --- Row 26: Ethereum Classic ---
Set-TextValue $ws.Range("D26") "26.91"
$ws.Range("E26").Value = "  +3.61%  "

# --- Row 27: Dai ---
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  -0.11%  "

# --- Row 28: Kaspa ---
$ws.Range("E28").Value = "  +0.28%  "

# --- Row 29: Cosmos ---
Set-TextValue $ws.Range("D29") "10.61"
$ws.Range("E29").Value = "  +2.93%  "

# --- Row 30: Injective Protocol ---
Set-TextValue $ws.Range("D30") "37.93"
$ws.Range("E30").Value = "  +3.42%  "

# --- Row 31: Toncoin ---
Set-TextValue $ws.Range("D31") "2.25"
$ws.Range("E31").Value = "  +10.10%  "

# --- Row 32: Filecoin ---
Set-TextValue $ws.Range("D32") "6.45"
$ws.Range("E32").Value = "  +4.42%  "

# --- Row 33: Render Token ---
Set-TextValue $ws.Range("D33") "6.13"
$ws.Range("E33").Value = "  +8.28%  "

# --- Row 34: OKB ---
Set-TextValue $ws.Range("D34") "53.10"
$ws.Range("E34").Value = "  +2.30%  "

# --- Row 35: Hedera ---
Set-TextValue $ws.Range("D35") "0.0936"
$ws.Range("E35").Value = "  +9.85%  "

# --- Row 36: VeChain ---
Set-TextValue $ws.Range("D36") "0.0452"
$ws.Range("E36").Value = "  +2.96%  "

# --- Row 37: First Digital USD ---
$ws.Range("E37").Value = "  -0.10%  "

# --- Row 38: Lido DAO Token ---
Set-TextValue $ws.Range("D38") "3.33"
$ws.Range("E38").Value = "  +6.55%  "

# --- Row 39: Celestia ---
Set-TextValue $ws.Range("D39") "18.80"
$ws.Range("E39").Value = "  +0.08%  "

# --- Row 40: ARBITRUM ---
$ws.Range("E40").Value = "  +4.84%  "

# --- Row 41: Stacks ---
$ws.Range("E41").Value = "  +10.76%  "

# --- Row 42: EnergySwap ---
Set-TextValue $ws.Range("D42") "24.02"
$ws.Range("E42").Value = "  +10.96%  "

# --- Row 43: Stellar ---
$ws.Range("E43").Value = "  +2.07%  "

# --- Row 44: Monero ---
Set-TextValue $ws.Range("D44") "122.65"
$ws.Range("E44").Value = "  +2.74%  "

# --- Row 45: WEMIX Token ---
Set-TextValue $ws.Range("D45") "2.20"
$ws.Range("E45").Value = "  +1.21%  "

# --- Row 46: NEAR Protocol ---
Set-TextValue $ws.Range("D46") "3.59"
$ws.Range("E46").Value = "  +6.16%  "

# --- Row 47: Maker ---
Set-TextValue $ws.Range("D47") "2.220.43"
$ws.Range("E47").Value = "  +4.80%  "

# --- Row 48: ApeX Protocol ---
Set-TextValue $ws.Range("D48") "2.53"
$ws.Range("E48").Value = "  +7.28%  "

# --- Row 49: The Graph ---
Set-TextValue $ws.Range("D49") "0.267"
$ws.Range("E49").Value = "  +24.94%  "

# --- Row 50: SEI ---
Set-TextValue $ws.Range("D50") "0.956"
$ws.Range("E50").Value = "  +5.68%  "

# --- Row 51: BEAM ---
Set-TextValue $ws.Range("D51") "0.0333"
$ws.Range("E51").Value = "  +14.63%  "
